$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy existing cell formats (that already carry the target style indices)
# onto the new cells, then set their values/text.

# K28 should match style of A2 (s="6")
$ws.Range("A2").Copy()
$ws.Range("K28").PasteSpecial(-4122)

# L28 should match style of P1 (s="16")
$ws.Range("P1").Copy()
$ws.Range("L28").PasteSpecial(-4122)

# K29 should match style of A3 (s="3")
$ws.Range("A3").Copy()
$ws.Range("K29").PasteSpecial(-4122)

# K30 should match style of A28 (s="14")
$ws.Range("A28").Copy()
$ws.Range("K30").PasteSpecial(-4122)

# K31 should match style of A5 (s="23")
$ws.Range("A5").Copy()
$ws.Range("K31").PasteSpecial(-4122)

# K32 should match style of K18 (s="29")
$ws.Range("K18").Copy()
$ws.Range("K32").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Now set the text values (L29:L32 keep the default/general style already)
$ws.Range("K28").Value = "TS 50-50 FS TA  "
$ws.Range("L28").Value = "All variables available with ExtraTreesClassifier applied"

$ws.Range("K29").Value = "TS 50-50 FS TA Top 7"
$ws.Range("L29").Value = "The top 7 variables from above test"

$ws.Range("K30").Value = "TS 50-50 TC"
$ws.Range("L30").Value = "No weather"

$ws.Range("K31").Value = "TS 50-50 TD"
$ws.Range("L31").Value = "No location info (besides Grid_Num)"

$ws.Range("K32").Value = "Logistic Regression"
$ws.Range("L32").Value = "Grid_Num, Hour, DayFrame, Join Count"

$ws.Range("N33").Select()
